$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new A1:D46 data in a 2D array (header row + 45 data rows)
$data = New-Object 'object[,]' 46,4

$data[0,0] = "time"
$data[0,1] = "LocationX"
$data[0,2] = "LocationY"
$data[0,3] = "LocationZ"
$data[1,0] = 0
$data[1,1] = 3361.16
$data[1,2] = 161.07
$data[1,3] = 3823.24
$data[2,0] = 0.03333333333333333
$data[2,1] = 4186.059999999999
$data[2,2] = 159.87
$data[2,3] = 4237.17
$data[3,0] = 0.06666666666666667
$data[3,1] = 4931.58
$data[3,2] = 187.57
$data[3,3] = 4619.700000000001
$data[4,0] = 0.1
$data[4,1] = 5584.19
$data[4,2] = 207.51
$data[4,3] = 4940.68
$data[5,0] = 0.1333333333333333
$data[5,1] = 6107.5
$data[5,2] = 225.18
$data[5,3] = 5180.700000000001
$data[6,0] = 0.1666666666666667
$data[6,1] = 6575.90673828125
$data[6,2] = 236.2182312011719
$data[6,3] = 5403.08056640625
$data[7,0] = 0.2
$data[7,1] = 6943.58837890625
$data[7,2] = 249.4012145996094
$data[7,3] = 5598.44287109375
$data[8,0] = 0.2333333333333333
$data[8,1] = 7354.48193359375
$data[8,2] = 267.7780151367188
$data[8,3] = 5918.787109375
$data[9,0] = 0.2666666666666667
$data[9,1] = 7820.8212890625
$data[9,2] = 274.03369140625
$data[9,3] = 6098.466796875
$data[10,0] = 0.3
$data[10,1] = 8161.99462890625
$data[10,2] = 294.5572814941406
$data[10,3] = 6165.677734375
$data[11,0] = 0.3333333333333333
$data[11,1] = 8260.9755859375
$data[11,2] = 300.0196533203125
$data[11,3] = 6194.69482421875
$data[12,0] = 0.3666666666666666
$data[12,1] = 8648.302734375
$data[12,2] = 316.169677734375
$data[12,3] = 6057.314453125
$data[13,0] = 0.4
$data[13,1] = 8980.2841796875
$data[13,2] = 310.091796875
$data[13,3] = 5855.31982421875
$data[14,0] = 0.4333333333333333
$data[14,1] = 9310.3203125
$data[14,2] = 328.6161499023438
$data[14,3] = 5954.09033203125
$data[15,0] = 0.4666666666666667
$data[15,1] = 9640.638671875
$data[15,2] = 333.1732177734375
$data[15,3] = 5976.60302734375
$data[16,0] = 0.5
$data[16,1] = 9579.0205078125
$data[16,2] = 347.687744140625
$data[16,3] = 6012.90576171875
$data[17,0] = 0.5333333333333333
$data[17,1] = 9763.2744140625
$data[17,2] = 369.9434204101562
$data[17,3] = 5961.96875
$data[18,0] = 0.5666666666666667
$data[18,1] = 10485.654296875
$data[18,2] = 371.8629150390625
$data[18,3] = 5948.7041015625
$data[19,0] = 0.6
$data[19,1] = 10408.322265625
$data[19,2] = 395.4696655273438
$data[19,3] = 6054.97900390625
$data[20,0] = 0.6333333333333333
$data[20,1] = 10924.951171875
$data[20,2] = 404.6325378417969
$data[20,3] = 6056.05517578125
$data[21,0] = 0.6666666666666666
$data[21,1] = 10754.6533203125
$data[21,2] = 432.5109252929688
$data[21,3] = 6068.1259765625
$data[22,0] = 0.7
$data[22,1] = 10950.1708984375
$data[22,2] = 429.2830505371094
$data[22,3] = 6007.95654296875
$data[23,0] = 0.7333333333333333
$data[23,1] = 11347.779296875
$data[23,2] = 444.1619873046875
$data[23,3] = 6088.14501953125
$data[24,0] = 0.7666666666666666
$data[24,1] = 11721.0654296875
$data[24,2] = 438.49658203125
$data[24,3] = 6001.56884765625
$data[25,0] = 0.8
$data[25,1] = 12283.880859375
$data[25,2] = 479.1244812011719
$data[25,3] = 6028.923828125
$data[26,0] = 0.8333333333333334
$data[26,1] = 12548.7451171875
$data[26,2] = 468.9779357910156
$data[26,3] = 5860.23388671875
$data[27,0] = 0.8666666666666667
$data[27,1] = 12832.2314453125
$data[27,2] = 490.4096984863281
$data[27,3] = 5889.4375
$data[28,0] = 0.9
$data[28,1] = 12980.322265625
$data[28,2] = 488.933837890625
$data[28,3] = 5626.94384765625
$data[29,0] = 0.9333333333333333
$data[29,1] = 13030.3251953125
$data[29,2] = 515.1332397460938
$data[29,3] = 5215.01953125
$data[30,0] = 0.9666666666666667
$data[30,1] = 13698.9150390625
$data[30,2] = 521.6781005859375
$data[30,3] = 5504.654296875
$data[31,0] = 1
$data[31,1] = 13692.474609375
$data[31,2] = 560.6834716796875
$data[31,3] = 5466.75341796875
$data[32,0] = 1.033333333333333
$data[32,1] = 13578.8896484375
$data[32,2] = 559.1596069335938
$data[32,3] = 4951.1943359375
$data[33,0] = 1.066666666666667
$data[33,1] = 14050.3671875
$data[33,2] = 575.5615234375
$data[33,3] = 4701.361328125
$data[34,0] = 1.1
$data[34,1] = 14201.75
$data[34,2] = 558.1874389648438
$data[34,3] = 4488.10009765625
$data[35,0] = 1.133333333333333
$data[35,1] = 14715.02734375
$data[35,2] = 600.772705078125
$data[35,3] = 4631.6962890625
$data[36,0] = 1.166666666666667
$data[36,1] = 14672.5234375
$data[36,2] = 601.9946899414062
$data[36,3] = 4328.28125
$data[37,0] = 1.2
$data[37,1] = 14801.24609375
$data[37,2] = 616.7559204101562
$data[37,3] = 4486.13525390625
$data[38,0] = 1.233333333333333
$data[38,1] = 15043.8623046875
$data[38,2] = 634.1895141601562
$data[38,3] = 3742.700439453125
$data[39,0] = 1.266666666666667
$data[39,1] = 15052.2470703125
$data[39,2] = 664.1375732421875
$data[39,3] = 3880.383056640625
$data[40,0] = 1.3
$data[40,1] = 15351.3115234375
$data[40,2] = 654.01220703125
$data[40,3] = 3647.98583984375
$data[41,0] = 1.333333333333333
$data[41,1] = 15331.12890625
$data[41,2] = 676.2843627929688
$data[41,3] = 3557.614990234375
$data[42,0] = 1.366666666666667
$data[42,1] = 15523.6337890625
$data[42,2] = 696.7689819335938
$data[42,3] = 3703.0693359375
$data[43,0] = 1.4
$data[43,1] = 15650.4189453125
$data[43,2] = 706.6704711914062
$data[43,3] = 3943.83837890625
$data[44,0] = 1.433333333333333
$data[44,1] = 15672.220703125
$data[44,2] = 721.937744140625
$data[44,3] = 4368.83349609375
$data[45,0] = 1.466666666666667
$data[45,1] = 15970.0244140625
$data[45,2] = 748.906494140625
$data[45,3] = 3908.42724609375

# Write header + first 45 data rows (rows 1-46) in one shot
$ws.Range("A1:D46").Value = $data

# Remove the old trailing row 47 (frame 45), which no longer exists in the new data
$ws.Rows.Item(47).Delete()

Write-Output "Updated header, refreshed A1:D46 values, removed row 47"
